$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 527735.6
$ws.Range("I33").Value = 1250204.9
$ws.Range("J33").Value = 2303.4546
$ws.Range("K33").Value = 1250204.9
$ws.Range("L33").Value = 2303.4546
$ws.Range("M33").Value = -1249975.9
$ws.Range("N33").Value = -2761.4546
$ws.Range("H40").Value = 3800
$ws.Range("I40").Value = 3175
$ws.Range("J40").Value = 4300
$ws.Range("K40").Value = 3175
$ws.Range("L40").Value = 4300
$ws.Range("M40").Value = -3000
$ws.Range("N40").Value = -4650
$ws.Range("H62").Value = 35029.08
$ws.Range("I62").Value = 44524.05
$ws.Range("K62").Value = 44524.05
$ws.Range("M62").Value = -43900.05
$ws.Range("H65").Value = 35029.08
$ws.Range("I65").Value = 44524.05
$ws.Range("K65").Value = 222620.25
$ws.Range("M65").Value = -219500.25
$ws.Range("H80").Value = 3385.3125
$ws.Range("J80").Value = 4139.4165
$ws.Range("L80").Value = 12418.2495
$ws.Range("N80").Value = -14414.2495
$ws.Range("H83").Value = 3385.3125
$ws.Range("J83").Value = 4139.4165
$ws.Range("L83").Value = 37254.7485
$ws.Range("N83").Value = -47238.7485
$ws.Range("H113").Value = 4599.1816
$ws.Range("I113").Value = 4002.25
$ws.Range("K113").Value = 4002.25
$ws.Range("M113").Value = -748.25
$ws.Range("H132").Value = 2457.682
$ws.Range("I132").Value = 2163.5
$ws.Range("K132").Value = 6490.5
$ws.Range("M132").Value = -3960.5
$ws.Range("H140").Value = 99725
$ws.Range("J140").Value = 99725
$ws.Range("L140").Value = 99725
$ws.Range("N140").Value = -110085

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4274.1113
$ws.Range("I74").Value = 4933.4
$ws.Range("K74").Value = 4933.4
$ws.Range("M74").Value = -4059.4
$ws.Range("H77").Value = 4274.1113
$ws.Range("I77").Value = 4933.4
$ws.Range("K77").Value = 24667
$ws.Range("M77").Value = -20299
$ws.Range("H102").Value = 8199.5
$ws.Range("I102").Value = 699
$ws.Range("J102").Value = 9699.6
$ws.Range("K102").Value = 699
$ws.Range("L102").Value = 9699.6
$ws.Range("M102").Value = 923
$ws.Range("N102").Value = -12943.6
$ws.Range("H122").Value = 2126.9333
$ws.Range("I122").Value = 1809.4546
$ws.Range("K122").Value = 5428.3638
$ws.Range("M122").Value = -2978.3638
$ws.Range("H140").Value = 63250
$ws.Range("J140").Value = 63250
$ws.Range("L140").Value = 63250
$ws.Range("N140").Value = -73610

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 39234.5
$ws.Range("H54").Value = 10681.4
$ws.Range("J54").Value = 45000
$ws.Range("L54").Value = 45000
$ws.Range("N54").Value = -45968
$ws.Range("H56").Value = 25000
$ws.Range("J56").Value = 25000
$ws.Range("L56").Value = 25000
$ws.Range("N56").Value = -26478
$ws.Range("H62").Value = 45714.43
$ws.Range("J62").Value = 45714.43
$ws.Range("L62").Value = 45714.43
$ws.Range("N62").Value = -47086.43
$ws.Range("H65").Value = 45714.43
$ws.Range("J65").Value = 45714.43
$ws.Range("L65").Value = 137143.29
$ws.Range("N65").Value = -144007.29
$ws.Range("H96").Value = 20426.834
$ws.Range("I96").Value = 14557
$ws.Range("K96").Value = 14557
$ws.Range("M96").Value = -11811
$ws.Range("H105").Value = 6904.522
$ws.Range("I105").Value = 7488.4707
$ws.Range("K105").Value = 7488.4707
$ws.Range("M105").Value = -5741.4707

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2979.2
$ws.Range("J22").Value = 2979.2
$ws.Range("L22").Value = 2979.2
$ws.Range("N22").Value = -3679.2
$ws.Range("H86").Value = 7240.615
$ws.Range("I86").Value = 7875.4287
$ws.Range("J86").Value = 6500
$ws.Range("K86").Value = 7875.4287
$ws.Range("L86").Value = 6500
$ws.Range("M86").Value = -6752.4287
$ws.Range("N86").Value = -8746
$ws.Range("H89").Value = 7240.615
$ws.Range("I89").Value = 7875.4287
$ws.Range("J89").Value = 6500
$ws.Range("K89").Value = 39377.14350000001
$ws.Range("L89").Value = 32500
$ws.Range("M89").Value = -33761.14350000001
$ws.Range("N89").Value = -43732

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 108
$ws.Range("I7").Value = 118.5
$ws.Range("J7").Value = 87
$ws.Range("K7").Value = 355.5
$ws.Range("L7").Value = 261
$ws.Range("M7").Value = -243.5
$ws.Range("N7").Value = -485
$ws.Range("H40").Value = 15.491228
$ws.Range("I40").Value = 14.205129
$ws.Range("J40").Value = 18.277779
$ws.Range("K40").Value = 56.820516
$ws.Range("L40").Value = 73.111116
$ws.Range("M40").Value = 12.179484
$ws.Range("N40").Value = -211.111116
$ws.Range("H101").Value = 10958.444
$ws.Range("J101").Value = 10958.444
$ws.Range("L101").Value = 32875.33199999999
$ws.Range("N101").Value = -37743.33199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 15181.091
$ws.Range("I20").Value = 9001
$ws.Range("J20").Value = 15799.1
$ws.Range("K20").Value = 9001
$ws.Range("L20").Value = 15799.1
$ws.Range("M20").Value = -8756
$ws.Range("N20").Value = -16289.1
$ws.Range("H80").Value = 6519.5
$ws.Range("J80").Value = 7531.7617
$ws.Range("L80").Value = 7531.7617
$ws.Range("N80").Value = -9527.761699999999
$ws.Range("H83").Value = 6519.5
$ws.Range("J83").Value = 7531.7617
$ws.Range("L83").Value = 37658.8085
$ws.Range("N83").Value = -47642.8085

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3210.375
$ws.Range("I82").Value = 228.33333
$ws.Range("K82").Value = 228.33333
$ws.Range("M82").Value = 132.66667
$ws.Range("H85").Value = 3210.375
$ws.Range("I85").Value = 228.33333
$ws.Range("K85").Value = 228.33333
$ws.Range("M85").Value = 1019.66667
$ws.Range("H100").Value = 6596.8335
$ws.Range("I100").Value = 5083.5293
$ws.Range("J100").Value = 10272
$ws.Range("K100").Value = 5083.5293
$ws.Range("L100").Value = 10272
$ws.Range("M100").Value = -4542.5293
$ws.Range("N100").Value = -11354
$ws.Range("H132").Value = 3534.15
$ws.Range("I132").Value = 2449.5454
$ws.Range("K132").Value = 7348.6362
$ws.Range("M132").Value = -4818.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1999.7858
$ws.Range("J81").Value = 3129.7144
$ws.Range("L81").Value = 6259.4288
$ws.Range("N81").Value = -8381.4288
$ws.Range("H84").Value = 1999.7858
$ws.Range("J84").Value = 3129.7144
$ws.Range("L84").Value = 31297.144
$ws.Range("N84").Value = -41905.144
$ws.Range("H100").Value = 3250
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -11082
$ws.Range("H107").Value = 534
$ws.Range("I107").Value = 290
$ws.Range("K107").Value = 870
$ws.Range("M107").Value = 1050
$ws.Range("H113").Value = 651.48
$ws.Range("I113").Value = 602.06665
$ws.Range("K113").Value = 1806.19995
$ws.Range("M113").Value = 363.8000500000001
$ws.Range("H133").Value = 50199.4
$ws.Range("J133").Value = 50199.4
$ws.Range("L133").Value = 50199.4
$ws.Range("N133").Value = -60319.4
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H141").Value = 200000
$ws.Range("I141").Value = 200000
$ws.Range("K141").Value = 200000
$ws.Range("M141").Value = -194820
